# Fruta / hortaliza, semanal
# Inserts a new weekly sample (date 2021-09-10 / serial 44449) for
# "Vega Monumental Concepción" - Plátano, as three quality rows
# (Maduro / Pintón / Primera Pintón), right above the existing block
# that starts at row 274. Everything below shifts down by 3 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows above row 274 (old data shifts down to 277+)
$ws.Range("A274:A276").EntireRow.Insert()

# Static (repeated) values shared by every row in this data block
$marketId   = 11
$market     = "Vega Monumental Concepción"
$region     = "Bíobío"
$codreg     = 8
$tipo       = "Fruta"
$prodId     = 100108
$producto   = "Tropicales y subtropicales"
$catId      = 100108006
$categoria  = "Plátano"
$variedad   = "Sin especificar"
$unidad     = "`$/caja 20 kilos"
$origen     = "Ecuador"
$kgUnidad   = 20

$fecha = 44449

$rows = @(
    @{ Row = 274; Calidad = "Maduro";         Volumen = 100; PMin = 13000; PMax = 13000; PProm = 13000; PKg = 650 },
    @{ Row = 275; Calidad = "Pintón";         Volumen = 300; PMin = 15000; PMax = 15000; PProm = 15000; PKg = 750 },
    @{ Row = 276; Calidad = "Primera Pintón"; Volumen = 300; PMin = 17000; PMax = 17000; PProm = 17000; PKg = 850 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $marketId
    $ws.Cells.Item($row, 2).Value  = $market
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $prodId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $catId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.PKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}
